$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.884.26"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "1.642.18"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.26"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.497"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0623"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "1.870.73"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.642.28"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.527"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.35"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").Value = "26.866.35"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "0.0₃0731"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.15"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.00"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.38"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.52"
$ws.Range("E22").Value = "  +4.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.39"
$ws.Range("E23").Value = "  -2.66%  "
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.68"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.23"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.76"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0509"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.20"
$ws.Range("E31").Value = "  +1.38%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("D35").Value = "1.272.26"
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.43"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.532"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.821"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.806"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.33"
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").Value = "1.781.79"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "92.65"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.06"
$ws.Range("E46").Value = "  -7.97%  "
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").Value = "0.0₆0101"
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.60"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0969"
$ws.Range("E51").Value = "  -0.08%  "
